$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Ligand/Receptor expression and derived edge-weight values
# recomputed from the new TPM input (NATMI Wnt4-Fzd2 output).
$updates = @{
    "G2" = 1.497411
    "H2" = 4.492233
    "I2" = 0.3090503153498128
    "J2" = 0.3090503153498128
    "M2" = 0.1811433333333334
    "N2" = 0.5434300000000001
    "O2" = 0.0111261749556462
    "P2" = 0.01112617495564619
    "Q2" = 0.27124601991
    "R2" = 2.44121417919
    "S2" = 0.003438547878679646
    "T2" = 0.003438547878679646
    "G3" = 1.497411
    "H3" = 4.492233
    "I3" = 0.3090503153498128
    "J3" = 0.3090503153498128
    "O3" = 0.8246098959508241
    "P3" = 0.8246098959508241
    "Q3" = 20.103238817178
    "R3" = 180.929149354602
    "S3" = 0.2548459483841785
    "T3" = 0.2548459483841785
    "G4" = 1.497411
    "H4" = 4.492233
    "I4" = 0.3090503153498128
    "J4" = 0.3090503153498128
    "M4" = 2.659118666666667
    "N4" = 7.977356
    "O4" = 0.1633282272592126
    "P4" = 0.1633282272592126
    "Q4" = 3.981793541771999
    "R4" = 35.836141875948
    "S4" = 0.05047664013998554
    "T4" = 0.05047664013998553
    "G5" = 1.497411
    "H5" = 4.492233
    "I5" = 0.3090503153498128
    "J5" = 0.3090503153498128
    "M5" = 0.015234
    "N5" = 0.045702
    "O5" = 0.0009357018343171013
    "P5" = 0.0009357018343171013
    "Q5" = 0.022811559174
    "R5" = 0.205304032566
    "S5" = 0.0002891789469690984
    "T5" = 0.0002891789469690984
    "I6" = 0.328071815935547
    "J6" = 0.3280718159355469
    "M6" = 0.1811433333333334
    "N6" = 0.5434300000000001
    "O6" = 0.0111261749556462
    "P6" = 0.01112617495564619
    "Q6" = 0.28794073294
    "R6" = 2.59146659646
    "S6" = 0.003650184422115451
    "T6" = 0.00365018442211545
    "I7" = 0.328071815935547
    "J7" = 0.3280718159355469
    "O7" = 0.8246098959508241
    "P7" = 0.8246098959508241
    "S7" = 0.2705312660030093
    "T7" = 0.2705312660030093
    "I8" = 0.328071815935547
    "J8" = 0.3280718159355469
    "M8" = 2.659118666666667
    "N8" = 7.977356
    "O8" = 0.1633282272592126
    "P8" = 0.1633282272592126
    "Q8" = 4.226865895447999
    "R8" = 38.041793059032
    "S8" = 0.05358338811046359
    "T8" = 0.05358338811046356
    "I9" = 0.328071815935547
    "J9" = 0.3280718159355469
    "M9" = 0.015234
    "N9" = 0.045702
    "O9" = 0.0009357018343171013
    "P9" = 0.0009357018343171013
    "Q9" = 0.024215570316
    "R9" = 0.217940132844
    "S9" = 0.0003069773999586338
    "T9" = 0.0003069773999586337
    "G10" = 1.151
    "H10" = 3.453
    "I10" = 0.2375546279329019
    "J10" = 0.2375546279329019
    "M10" = 0.1811433333333334
    "N10" = 0.5434300000000001
    "O10" = 0.0111261749556462
    "P10" = 0.01112617495564619
    "Q10" = 0.2084959766666667
    "R10" = 1.87646379
    "S10" = 0.002643074351904903
    "T10" = 0.002643074351904903
    "G11" = 1.151
    "H11" = 3.453
    "I11" = 0.2375546279329019
    "J11" = 0.2375546279329019
    "O11" = 0.8246098959508241
    "P11" = 0.8246098959508241
    "Q11" = 15.45255636466667
    "R11" = 139.073007282
    "S11" = 0.195889897022387
    "T11" = 0.1958898970223869
    "G12" = 1.151
    "H12" = 3.453
    "I12" = 0.2375546279329019
    "J12" = 0.2375546279329019
    "M12" = 2.659118666666667
    "N12" = 7.977356
    "O12" = 0.1633282272592126
    "P12" = 0.1633282272592126
    "Q12" = 3.060645585333333
    "R12" = 27.545810268
    "S12" = 0.0387993762575027
    "T12" = 0.03879937625750268
    "G13" = 1.151
    "H13" = 3.453
    "I13" = 0.2375546279329019
    "J13" = 0.2375546279329019
    "M13" = 0.015234
    "N13" = 0.045702
    "O13" = 0.0009357018343171013
    "P13" = 0.0009357018343171013
    "Q13" = 0.017534334
    "R13" = 0.157809006
    "S13" = 0.0002222803011073328
    "T13" = 0.0002222803011073328
    "G14" = 0.6072163333333332
    "H14" = 1.821649
    "I14" = 0.1253232407817384
    "J14" = 0.1253232407817384
    "M14" = 0.1811433333333334
    "N14" = 0.5434300000000001
    "O14" = 0.0111261749556462
    "P14" = 0.01112617495564619
    "Q14" = 0.1099931906744444
    "R14" = 0.9899387160700001
    "S14" = 0.001394368302946196
    "T14" = 0.001394368302946196
    "G15" = 0.6072163333333332
    "H15" = 1.821649
    "I15" = 0.1253232407817384
    "J15" = 0.1253232407817384
    "O15" = 0.8246098959508241
    "P15" = 0.8246098959508241
    "Q15" = 8.152080466011776
    "R15" = 73.36872419410599
    "S15" = 0.1033427845412494
    "T15" = 0.1033427845412494
    "G16" = 0.6072163333333332
    "H16" = 1.821649
    "I16" = 0.1253232407817384
    "J16" = 0.1253232407817384
    "M16" = 2.659118666666667
    "N16" = 7.977356
    "O16" = 0.1633282272592126
    "P16" = 0.1633282272592126
    "Q16" = 1.614660286671555
    "R16" = 14.531942580044
    "S16" = 0.02046882275126079
    "T16" = 0.02046882275126079
    "G17" = 0.6072163333333332
    "H17" = 1.821649
    "I17" = 0.1253232407817384
    "J17" = 0.1253232407817384
    "M17" = 0.015234
    "N17" = 0.045702
    "O17" = 0.0009357018343171013
    "P17" = 0.0009357018343171013
    "Q17" = 0.009250333621999998
    "R17" = 0.08325300259799999
    "S17" = 0.0001172651862820364
    "T17" = 0.0001172651862820364
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
